$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row right after the header row (row 1), shifting existing
# data rows (old rows 2-4) down by one. This new row carries the
# "slug"/identifier form of each column's header, used to relate columns
# into SKOS hierarchies.
$ws.Rows.Item(2).Insert()

$ws.Cells.Item(2, 1).Value = "n-hogares"
$ws.Cells.Item(2, 2).Value = "tipo-de-hogar-codigo"
$ws.Cells.Item(2, 3).Value = "municipio-codigo"
$ws.Cells.Item(2, 4).Value = "tipo-de-hogar"
$ws.Cells.Item(2, 5).Value = "municipio-nombre"
